# Update the Handoff/Handback datetime stamps for the first data row
# (row 2) on the "zh-cn" and "de-de" worksheets, reflecting a newly
# generated handback report.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-19 18:52:22"
$wsZhCn.Range("H2").Value = "2016-03-19 18:52:41"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-19 18:52:26"
$wsDeDe.Range("H2").Value = "2016-03-19 18:52:46"
